$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form responses 1")
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------
# Add two new rows to the "Form_Responses" table (auto-expands the table
# range / worksheet dimension, like typing new rows below an Excel Table).
# ---------------------------------------------------------------------
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# ---------------------------------------------------------------------
# Formatting: mirror the alternating banding pattern already used by the
# table. Row 23 is now the last row, so it takes over the special
# "final row" bottom-border look that row 21 used to have; row 21 becomes
# a normal "odd" row (like row 19); row 22 is a normal "even" row
# (like row 20). Column I (LinkedIn hyperlink) and J (PMI ID) are handled
# separately since rows 22/23 don't use both of them identically.
# ---------------------------------------------------------------------
$ws.Range("A21:H21").Copy()
$ws.Range("A23:H23").PasteSpecial(-4122)
$ws.Range("K21:AD21").Copy()
$ws.Range("K23:AD23").PasteSpecial(-4122)

$ws.Range("A20:I20").Copy()
$ws.Range("A22:I22").PasteSpecial(-4122)
$ws.Range("K20:AD20").Copy()
$ws.Range("K22:AD22").PasteSpecial(-4122)

$ws.Range("A19:AD19").Copy()
$ws.Range("A21:AD21").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 22: Jacqueline Shen
# ---------------------------------------------------------------------
$ws.Range("A22").Value = 45922.68002349537
$ws.Range("B22").Value = "jacqueline.shen@blueyonder.com"
$ws.Range("C22").Value = "Jacqueline"
$ws.Range("D22").Value = "Shen"
$ws.Range("E22").Value = "'0424351713"
$ws.Range("F22").Value = "jacqueline.shen@blueyonder.com"
$ws.Range("G22").Value = "marketing coordinator"
$ws.Range("H22").Value = "Blue Yonder"
$ws.Range("K22").Value = "Software Implementation and/or Rollout, Strategic Planning, Business Change Management or Change Related Initiatives, Fund Raising Initiatives, Problem Analysis and Solutioning, Expansion of Membership and/or Increasing Awareness and Support for NFP Organisation, Events Planning"
$ws.Range("L22").Value = "Yes"
$ws.Range("M22").Value = "None"
$ws.Range("N22").Value = "N/A"
$ws.Range("O22").Value = "4 - 8 Years"
$ws.Range("P22").Value = 4
$ws.Range("Q22").Value = 4
$ws.Range("R22").Value = 2
$ws.Range("S22").Value = 3
$ws.Range("T22").Value = 4
$ws.Range("U22").Value = 2
$ws.Range("V22").Value = 2
$ws.Range("W22").Value = 2
$ws.Range("X22").Value = 3
$ws.Range("Y22").Value = 2
$ws.Range("Z22").Value = 3
$ws.Range("AA22").Value = 5
$ws.Range("AB22").Value = 2
$ws.Range("AC22").Value = "Yes"
$ws.Range("AD22").Value = "Jacqueline Shen"

# LinkedIn hyperlink for row 22 (column I)
$ws.Hyperlinks.Add($ws.Range("I22"), "https://www.linkedin.com/in/jacqueline-shen-5b565917a/") | Out-Null
$ws.Range("I20").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 23: Nikki Gittins
# ---------------------------------------------------------------------
$ws.Range("A23").Value = 45922.751129837961
$ws.Range("B23").Value = "ngittins@atlassian.com"
$ws.Range("C23").Value = "Nikki"
$ws.Range("D23").Value = "Gittins"
$ws.Range("E23").Value = "'0409198268"
$ws.Range("F23").Value = "ngittins@atlassian.com"
$ws.Range("G23").Value = "Snr Principle Program Manager [Technical] "
$ws.Range("H23").Value = "Atlassian"
$ws.Range("K23").Value = "Software Implementation and/or Rollout, Strategic Planning, Business Change Management or Change Related Initiatives, Problem Analysis and Solutioning, Events Planning"
$ws.Range("L23").Value = "Yes"
$ws.Range("M23").Value = "None"
$ws.Range("N23").Value = "None"
$ws.Range("O23").Value = "More than 8 Years"
$ws.Range("P23").Value = 5
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 5
$ws.Range("S23").Value = 5
$ws.Range("T23").Value = 5
$ws.Range("U23").Value = 5
$ws.Range("V23").Value = 5
$ws.Range("W23").Value = 5
$ws.Range("X23").Value = 5
$ws.Range("Y23").Value = 5
$ws.Range("Z23").Value = 1
$ws.Range("AA23").Value = 5
$ws.Range("AB23").Value = 5
$ws.Range("AC23").Value = "Yes"
$ws.Range("AD23").Value = "Nikki Gittins"

# ---------------------------------------------------------------------
# Misc cosmetic bits from the diff
# ---------------------------------------------------------------------
$ws.Range("G6").Select() | Out-Null

Write-Host "done"
